# Update the "Förändrad" (Changed) date column C for rows 2-24
# from serial date 45224 to 45225 (one day later).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 24; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45224) {
        $cell.Value = 45225
    }
}
